# PETAPP_SBL_02.xlsx — "Agregacion de el sprint burndown 1 y 2, parte grafica"
#
# This script takes the "SPRINT BACKLOG 1" sheet (which actually holds the
# content for "Sprint Backlog 2") and:
#   - clears out the now-unused Day 4..Day 7 columns (F:I) in the header and
#     the data rows, leaving only Day 1..Day 3 (C:E)
#   - corrects a couple of hour values in the remaining Day 1..Day 3 columns
#   - adds a totals row ("Nro Horas") under the task rows with SUM formulas
#   - adds a second author
#   - adds a Sprint Burndown line chart driven by the new totals row
#   - tidies up the selection / view state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPRINT BACKLOG 1")

# ---------------------------------------------------------------------
# Row 7 (logout client) and Row 9 (functional tests) hour corrections
# ---------------------------------------------------------------------
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("E9").Value = 0

# ---------------------------------------------------------------------
# Clear the Day 4..Day 7 columns (F:I) — header row and data rows
# ---------------------------------------------------------------------
$ws.Range("F4:I4").ClearContents()
$ws.Range("F6:I10").ClearContents()

# ---------------------------------------------------------------------
# Totals row — "Nro Horas" with SUM formulas over the task rows
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "Nro Horas"
$ws.Range("C11").Formula = "=SUM(C6:C10)"
$ws.Range("D11:E11").Formula = "=SUM(D6:D10)"

# ---------------------------------------------------------------------
# Second author
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "Carlos Zarate Carpio"

# ---------------------------------------------------------------------
# Sprint Burndown 2 chart
# ---------------------------------------------------------------------
$chartObjects = $ws.ChartObjects()
$co = $chartObjects.Add(66675, 1143000, 309232, 900430)
$chart = $co.Chart
$chart.ChartType = 4
$chart.SetSourceData($ws.Range("C11:E11"))

$series = $chart.SeriesCollection(1)
$series.Name = "='SPRINT BACKLOG 1'!`$B`$11"
$series.XValues = $ws.Range("C4:E4")

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Sprint Burndown 2"
$chart.HasLegend = $false

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "Nro Dias"
$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "Nro Horas"

# ---------------------------------------------------------------------
# View / selection tidy-up
# ---------------------------------------------------------------------
$ws.Range("I26").Select()
